# Document Sprint 8 #65
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Copy the Sprint 7 block (title rows + table, rows 99:110) down to rows 113:124
$src = $ws.Range("B99:G110")
$src.Copy()
$dst = $ws.Range("B113")
$dst.PasteSpecial()

# Update the title text for the new block
$ws.Range("B113").Value = "Sprint 8"

# Update the "Min." (D) and "Real" (F) columns with the actual sprint-8 data
$ws.Range("D116").Value = 0
$ws.Range("D117").Value = 0
$ws.Range("D118").Value = 0
$ws.Range("D119").Value = 2
$ws.Range("D120").Value = 3
$ws.Range("D121").Value = 0
$ws.Range("D122").Value = 0
$ws.Range("D123").Value = 0

$ws.Range("F116").Value = 0
$ws.Range("F117").Value = 0
$ws.Range("F118").Value = 0
$ws.Range("F119").Value = 1
$ws.Range("F120").Value = 4
$ws.Range("F121").Value = 0
$ws.Range("F122").Value = 0
$ws.Range("F123").Value = 0

$ws.Range("F128").Select()
